$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Members (row 5) and People (row 6) CRUD columns: Edit/View/Delete now pass ---
$ws.Range("D5:G6").Value = 1

# --- Column B gets a pass/fail/partial rollup formula based on the CRUD checks in C:H ---
# Row 3 (Bills): standalone formula
$ws.Range("B3").Formula = "=IF(SUM(C3:H3)=6,1,IF(SUM(C3:H3)>-6,0,-1))"

# Rows 4-22: same formula pattern, one per row (source workbook stores these as one
# shared-formula family anchored at B4)
for ($r = 4; $r -le 22; $r++) {
    $ws.Range("B$r").Formula = "=IF(SUM(C$r" + ":H$r" + ")=6,1,IF(SUM(C$r" + ":H$r" + ")>-6,0,-1))"
}

# --- Conditional formatting: column B now participates in the icon-set rule too, so
#     the banding is re-sliced (B3:H3 / C4:H15 / B4:B22) and the Courses row (16) loses
#     its own column-B cell from the lower band (C16:H22) ---
$fcTop = $ws.Range("B3:H15").FormatConditions.Item(1)
$fcTop.ModifyAppliesToRange($ws.Range("B3:H3"))

$fcMid = $ws.Range("C4:H15").FormatConditions.AddIconSetCondition()
$fcMid.IconSet = "3Symbols"

$fcCol = $ws.Range("B4:B22").FormatConditions.AddIconSetCondition()
$fcCol.IconSet = "3Symbols"

$fcLow = $ws.Range("B16:H22").FormatConditions.Item(1)
$fcLow.ModifyAppliesToRange($ws.Range("C16:H22"))

# --- Selection moved to D6 ---
$ws.Range("D6").Select()
